# Apply release-note-summary.xlsx fixes on the "uat" worksheet:
#  - F3: remove stray literal double-quotes around AQST -> AQST
#  - F4/F5/F6: fix JSON key typo "repouatory" -> "repository"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("uat")

$ws.Range("F3").Value = "AQST"

$f4 = $ws.Range("F4").Value()
$f4 = $f4.Replace("repouatory", "repository")
$ws.Range("F4").Value = $f4

$f5 = $ws.Range("F5").Value()
$f5 = $f5.Replace("repouatory", "repository")
$ws.Range("F5").Value = $f5

$f6 = $ws.Range("F6").Value()
$f6 = $f6.Replace("repouatory", "repository")
$ws.Range("F6").Value = $f6
